$wb = $excel.ActiveWorkbook

# Insert the new "Hoja2" worksheet right after "reto2" (so the order becomes
# Hoja1, reto2, Hoja2, reemp) and make it the active sheet/tab, exactly like
# Excel does when a user inserts a new sheet via the UI.
$after = $wb.Worksheets.Item("reto2")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$ws.Name = "Hoja2"

# Row 1
$ws.Range("A1").Value = "ACT"
$ws.Range("B1").Value = "7 Helado 65000 11"
$ws.Range("C1").Value = "Helado Galletas 10950.0 1544600.0"
$ws.Range("E1").Value = "ERROR"
$ws.Range("G1").Formula = '=+A1&"  "&B1&"\t"&C1&"\t"&E1&"  "&F1'

# Row 2
$ws.Range("A2").Value = "ACT"
$ws.Range("B2").Value = "10 Jamon 500 10 "
$ws.Range("C2").Value = "Arandanos Galletas 3450.0 869100.0"
$ws.Range("E2").Value = "ACT"
$ws.Range("F2").Value = "15 Papa 1500 10 "
$ws.Range("G2").Formula = '=+A2&"  "&B2&"\t"&C2&"\t"&E2&"  "&F2'

# Row 3
$ws.Range("A3").Value = "BOR"
$ws.Range("B3").Value = "10 Jamon 15000 10"
$ws.Range("C3").Value = "Arandanos Galletas 3777.8 864100.0"
$ws.Range("E3").Value = "AGR"
$ws.Range("F3").Value = "3 Peras 2700 33 "
$ws.Range("G3").Formula = '=+A3&"  "&B3&"\t"&C3&"\t"&E3&"  "&F3'

# Row 4
$ws.Range("A4").Value = "BOR"
$ws.Range("B4").Value = "3 Peras 2700 33 "
$ws.Range("C4").Value = "Jamon Galletas 5144.4 925000.0"
$ws.Range("E4").Value = "BOR"
$ws.Range("F4").Value = "15 Papa 1500 10"
$ws.Range("G4").Formula = '=+A4&"  "&B4&"\t"&C4&"\t\t"&E4&"  "&F4'

# Row 5
$ws.Range("A5").Value = "AGR"
$ws.Range("B5").Value = "11 Melon 70 13"
$ws.Range("C5").Value = "Jamon Melon 4460.9 1015010.0"
$ws.Range("E5").Value = "BOR"
$ws.Range("F5").Value = "14 Maiz 45000 12 "
$ws.Range("G5").Formula = '=+A5&"  "&B5&"\t"&C5&"\t\t"&E5&"  "&F5'

# Row 6
$ws.Range("A6").Value = "AGR"
$ws.Range("B6").Value = "11 Maiz 70000 1 "
$ws.Range("C6").Value = "Maiz Galletas 10818.2 1084100.0"
$ws.Range("G6").Formula = '=+A6&"  "&B6&"\t"&C6&"\t\t"&E6&"  "&F6'

# Leave the selection on G1, matching the saved cursor position.
[void]$ws.Range("G1").Select()
